$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.068598123028721
$ws.Range("D2").Value = 1.069214107346037
$ws.Range("E2").Value = 1.072716216237676
$ws.Range("F2").Value = 1.081869338050028
$ws.Range("I2").Value = 1.053098109137426
$ws.Range("J2").Value = 1.073537116518704
$ws.Range("K2").Value = 1.07191721343545
$ws.Range("L2").Value = 1.075409992178965
$ws.Range("M2").Value = 1.084539036760271

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.069752685081886
$ws.Range("D3").Value = 1.070115331736703
$ws.Range("E3").Value = 1.073732100620997
$ws.Range("F3").Value = 1.082935817633667
$ws.Range("I3").Value = 1.053430970584871
$ws.Range("J3").Value = 1.07434753643135
$ws.Range("K3").Value = 1.072634353799943
$ws.Range("L3").Value = 1.076242183385786
$ws.Range("M3").Value = 1.085423441781775

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.070499949551977
$ws.Range("D4").Value = 1.070698609005516
$ws.Range("E4").Value = 1.074389896539057
$ws.Range("F4").Value = 1.083626430059325
$ws.Range("I4").Value = 1.053645243351831
$ws.Range("J4").Value = 1.074871529106607
$ws.Range("K4").Value = 1.073097879451309
$ws.Range("L4").Value = 1.076780502638638
$ws.Range("M4").Value = 1.085995633006102

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.070814144591256
$ws.Range("D5").Value = 1.070943848465059
$ws.Range("E5").Value = 1.074666541792543
$ws.Range("F5").Value = 1.083916890066405
$ws.Range("I5").Value = 1.053735057681414
$ws.Range("J5").Value = 1.075091719204988
$ws.Range("K5").Value = 1.073292623020768
$ws.Range("L5").Value = 1.077006772778794
$ws.Range("M5").Value = 1.08623616352328

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.070866901981096
$ws.Range("D6").Value = 1.070985026985379
$ws.Range("E6").Value = 1.074712998070245
$ws.Range("F6").Value = 1.083965666998632
$ws.Range("I6").Value = 1.053750122307585
$ws.Range("J6").Value = 1.075128684440242
$ws.Range("K6").Value = 1.073325314132524
$ws.Range("L6").Value = 1.077044762228376
$ws.Range("M6").Value = 1.086276548560273

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.070504147663637
$ws.Range("D7").Value = 1.070701885790722
$ws.Range("E7").Value = 1.07439359266278
$ws.Range("F7").Value = 1.08363031070299
$ws.Range("I7").Value = 1.053646444500264
$ws.Range("J7").Value = 1.074874471676213
$ws.Range("K7").Value = 1.073100482106216
$ws.Range("L7").Value = 1.07678352622562
$ws.Range("M7").Value = 1.085998847060023

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.068988274617785
$ws.Range("D8").Value = 1.06951865438978
$ws.Range("E8").Value = 1.07305944573912
$ws.Range("F8").Value = 1.082229650287652
$ws.Range("I8").Value = 1.053210831228347
$ws.Range("J8").Value = 1.07381108527012
$ws.Range("K8").Value = 1.072159680409149
$ws.Range("L8").Value = 1.0756912690804
$ws.Range("M8").Value = 1.08483794171514

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.06631850633067
$ws.Range("D9").Value = 1.067434613358311
$ws.Range("E9").Value = 1.070711963607996
$ws.Range("F9").Value = 1.079765559901798
$ws.Range("I9").Value = 1.052434717154804
$ws.Range("J9").Value = 1.071934169495964
$ws.Range("K9").Value = 1.070497945570388
$ws.Range("L9").Value = 1.073765314589347
$ws.Range("M9").Value = 1.082791681043411

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.064539556177625
$ws.Range("D10").Value = 1.066045908775603
$ws.Range("E10").Value = 1.069149299354002
$ws.Range("F10").Value = 1.078125557279785
$ws.Range("I10").Value = 1.051911586127511
$ws.Range("J10").Value = 1.070680803198111
$ws.Range("K10").Value = 1.069387483238042
$ws.Range("L10").Value = 1.07248049525654
$ws.Range("M10").Value = 1.081427103033606

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.063769453377602
$ws.Range("D11").Value = 1.065444740842335
$ws.Range("E11").Value = 1.068473198884212
$ws.Range("F11").Value = 1.077416063191964
$ws.Range("I11").Value = 1.05168370674739
$ws.Range("J11").Value = 1.070137583336885
$ws.Range("K11").Value = 1.068906013591176
$ws.Range("L11").Value = 1.071923951071887
$ws.Range("M11").Value = 1.080836128338912

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.063483431132369
$ws.Range("D12").Value = 1.065221462857187
$ws.Range("E12").Value = 1.068222146543281
$ws.Range("F12").Value = 1.077152621134449
$ws.Range("I12").Value = 1.051598857713459
$ws.Range("J12").Value = 1.069935731226748
$ws.Range("K12").Value = 1.06872707915533
$ws.Range("L12").Value = 1.071717194147251
$ws.Range("M12").Value = 1.080616598258584

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.063544782613359
$ws.Range("D13").Value = 1.065269355722892
$ws.Range("E13").Value = 1.068275994438386
$ws.Range("F13").Value = 1.077209126025422
$ws.Range("I13").Value = 1.051617067374778
$ws.Range("J13").Value = 1.069979032637344
$ws.Range("K13").Value = 1.068765465509158
$ws.Range("L13").Value = 1.071761545644865
$ws.Range("M13").Value = 1.080663688923117

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.063745810119465
$ws.Range("D14").Value = 1.065426284154547
$ws.Range("E14").Value = 1.068452445158282
$ws.Range("F14").Value = 1.077394285039472
$ws.Range("I14").Value = 1.051676697271981
$ws.Range("J14").Value = 1.070120899732243
$ws.Range("K14").Value = 1.068891224757519
$ws.Range("L14").Value = 1.071906861126735
$ws.Range("M14").Value = 1.080817982236922

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.063869673511489
$ws.Range("D15").Value = 1.065522975928064
$ws.Range("E15").Value = 1.068561173061767
$ws.Range("F15").Value = 1.077508380286038
$ws.Range("I15").Value = 1.051713410119414
$ws.Range("J15").Value = 1.070208298641625
$ws.Range("K15").Value = 1.068968696563867
$ws.Range("L15").Value = 1.071996390592003
$ws.Range("M15").Value = 1.080913045343965

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.064590669375619
$ws.Range("D16").Value = 1.066085809475758
$ws.Range("E16").Value = 1.069194181350923
$ws.Range("F16").Value = 1.078172657488174
$ws.Range("I16").Value = 1.051926681062425
$ws.Range("J16").Value = 1.070716844263062
$ws.Range("K16").Value = 1.069419423457371
$ws.Range("L16").Value = 1.072517426838242
$ws.Range("M16").Value = 1.081466321879406

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.06504298237719
$ws.Range("D17").Value = 1.066438900312277
$ws.Range("E17").Value = 1.069591396105844
$ws.Range("F17").Value = 1.078589511664597
$ws.Range("I17").Value = 1.0520600960342
$ws.Range("J17").Value = 1.071035706318366
$ws.Range("K17").Value = 1.069701983148324
$ws.Range("L17").Value = 1.072844202999432
$ws.Range("M17").Value = 1.081813349784479

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.065306827716249
$ws.Range("D18").Value = 1.066644866728696
$ws.Range("E18").Value = 1.069823137199998
$ws.Range("F18").Value = 1.078832717211375
$ws.Range("I18").Value = 1.052137783499864
$ws.Range("J18").Value = 1.071221644556738
$ws.Range("K18").Value = 1.069866734457671
$ws.Range("L18").Value = 1.07303478586275
$ws.Range("M18").Value = 1.082015755327843

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.065396795380062
$ws.Range("D19").Value = 1.066715098415334
$ws.Range("E19").Value = 1.069902163825714
$ws.Range("F19").Value = 1.078915654465153
$ws.Range("I19").Value = 1.052164250663171
$ws.Range("J19").Value = 1.071285036457551
$ws.Range("K19").Value = 1.069922900059846
$ws.Range("L19").Value = 1.073099766284903
$ws.Range("M19").Value = 1.082084768704692

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.064994451533234
$ws.Range("D20").Value = 1.0664010155069
$ws.Range("E20").Value = 1.06954877329803
$ws.Range("F20").Value = 1.07854478079056
$ws.Range("I20").Value = 1.052045795445005
$ws.Range("J20").Value = 1.071001500467231
$ws.Range("K20").Value = 1.06967167347618
$ws.Range("L20").Value = 1.072809145095039
$ws.Range("M20").Value = 1.081776118006501

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.063686611758109
$ws.Range("D21").Value = 1.065380072020457
$ws.Range("E21").Value = 1.068400482571409
$ws.Range("F21").Value = 1.077339757686321
$ws.Range("I21").Value = 1.051659143399351
$ws.Range("J21").Value = 1.070079125511349
$ws.Range("K21").Value = 1.068854194422154
$ws.Range("L21").Value = 1.071864070214434
$ws.Range("M21").Value = 1.080772547138706

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.062864483795812
$ws.Range("D22").Value = 1.064738294794951
$ws.Range("E22").Value = 1.067678977931866
$ws.Range("F22").Value = 1.076582665063283
$ws.Range("I22").Value = 1.051414856756399
$ws.Range("J22").Value = 1.069498751497278
$ws.Range("K22").Value = 1.068339662163455
$ws.Range("L22").Value = 1.071269680980905
$ws.Range("M22").Value = 1.080141470975535

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.063300294276089
$ws.Range("D23").Value = 1.065078500688652
$ws.Range("E23").Value = 1.068061416574682
$ws.Range("F23").Value = 1.076983961812737
$ws.Range("I23").Value = 1.051544469920947
$ws.Range("J23").Value = 1.069806460591793
$ws.Range("K23").Value = 1.068612477718944
$ws.Range("L23").Value = 1.071584795401161
$ws.Range("M23").Value = 1.080476025122901

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.06501638047416
$ws.Range("D24").Value = 1.066418133974767
$ws.Range("E24").Value = 1.069568032547984
$ws.Range("F24").Value = 1.078564992556742
$ws.Range("I24").Value = 1.052052257671087
$ws.Range("J24").Value = 1.071016956768975
$ws.Range("K24").Value = 1.069685369300617
$ws.Range("L24").Value = 1.072824986315055
$ws.Range("M24").Value = 1.081792941474959

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.067008544571607
$ws.Range("D25").Value = 1.067973272865024
$ws.Range("E25").Value = 1.071318434472575
$ws.Range("F25").Value = 1.080402106142701
$ws.Range("I25").Value = 1.052636369625016
$ws.Range("J25").Value = 1.072419764741992
$ws.Range("K25").Value = 1.070928008864739
$ws.Range("L25").Value = 1.074263369463477
$ws.Range("M25").Value = 1.083320759525489
